# Trading update: 2026-02-18 10:30:00
#
# New trades opened since the last export:
#   #22 MarketMaking DOWN 0.12        10:28:30
#   #23 MarketMaking UP   0.92        10:29:21
#   #24 MarketMaking DOWN 0.04        10:29:26
#   #25 MarketMaking DOWN 0.388544    10:29:32
#   #26 momentum     UP   0.42        10:29:39
#
# "All Trades" gets every new row appended; the per-strategy detail
# columns (Capital After / Slippage / Confidence / Entry Reason /
# Duration) only stay populated for the trades that are still being
# actively tracked on the strategy's own sheet (last 3 for
# MarketMaking, last 1 for momentum) - older rows for a strategy lose
# those values once they roll off that sheet.  The "momentum" and
# "MarketMaking" sheets themselves are overwritten in place with the
# latest open trades for that strategy.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, $row, $col, [string]$text)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------
# "All Trades" sheet
# ---------------------------------------------------------------
$all = $wb.Worksheets.Item("All Trades")

# Rows 19-22 (trades #18-#21) age out of the rolling per-strategy
# window: Exit Price becomes an explicit 0 and the capital/slippage/
# confidence/reason/duration columns go blank again.
foreach ($r in 19, 20, 21, 22) {
    $all.Cells.Item($r, 7).Value = 0          # G: Exit Price
    $all.Range("K" + $r + ":O" + $r).ClearContents()
    $all.Range("Q" + $r).ClearContents()
}

# New trade #22 -> row 23 (MarketMaking, still within the rolling
# window boundary, so it is freshly OPEN with no derived columns yet)
Set-TextCell $all 23 2 "2026-02-18"
Set-TextCell $all 23 3 "10:28:30"
$all.Cells.Item(23, 1).Value = 22
$all.Cells.Item(23, 4).Value = "MarketMaking"
$all.Cells.Item(23, 5).Value = "DOWN"
$all.Cells.Item(23, 6).Value = 0.12
$all.Cells.Item(23, 7).Value = 0
$all.Cells.Item(23, 8).Value = "OPEN"
$all.Cells.Item(23, 9).Value = 0
$all.Cells.Item(23, 10).Value = 0

# New trade #23 -> row 24 (MarketMaking)
Set-TextCell $all 24 2 "2026-02-18"
Set-TextCell $all 24 3 "10:29:21"
$all.Cells.Item(24, 1).Value = 23
$all.Cells.Item(24, 4).Value = "MarketMaking"
$all.Cells.Item(24, 5).Value = "UP"
$all.Cells.Item(24, 6).Value = 0.92
$all.Cells.Item(24, 8).Value = "OPEN"
$all.Cells.Item(24, 9).Value = 0
$all.Cells.Item(24, 10).Value = 0
$all.Cells.Item(24, 11).Value = 100
$all.Cells.Item(24, 12).Value = 0
$all.Cells.Item(24, 13).Value = 0
$all.Cells.Item(24, 14).Value = 0.6
$all.Cells.Item(24, 15).Value = "Normal spread capture: 408 bps"
$all.Cells.Item(24, 17).Value = 0

# New trade #24 -> row 25 (MarketMaking)
Set-TextCell $all 25 2 "2026-02-18"
Set-TextCell $all 25 3 "10:29:26"
$all.Cells.Item(25, 1).Value = 24
$all.Cells.Item(25, 4).Value = "MarketMaking"
$all.Cells.Item(25, 5).Value = "DOWN"
$all.Cells.Item(25, 6).Value = 0.04
$all.Cells.Item(25, 8).Value = "OPEN"
$all.Cells.Item(25, 9).Value = 0
$all.Cells.Item(25, 10).Value = 0
$all.Cells.Item(25, 11).Value = 100
$all.Cells.Item(25, 12).Value = 0
$all.Cells.Item(25, 13).Value = 0
$all.Cells.Item(25, 14).Value = 0.6
$all.Cells.Item(25, 15).Value = "Normal spread capture: 202 bps"
$all.Cells.Item(25, 17).Value = 0

# New trade #25 -> row 26 (MarketMaking)
Set-TextCell $all 26 2 "2026-02-18"
Set-TextCell $all 26 3 "10:29:32"
$all.Cells.Item(26, 1).Value = 25
$all.Cells.Item(26, 4).Value = "MarketMaking"
$all.Cells.Item(26, 5).Value = "DOWN"
$all.Cells.Item(26, 6).Value = 0.388544
$all.Cells.Item(26, 8).Value = "OPEN"
$all.Cells.Item(26, 9).Value = 0
$all.Cells.Item(26, 10).Value = 0
$all.Cells.Item(26, 11).Value = 100
$all.Cells.Item(26, 12).Value = 0
$all.Cells.Item(26, 13).Value = 0
$all.Cells.Item(26, 14).Value = 0.6
$all.Cells.Item(26, 15).Value = "Normal spread capture: 202 bps"
$all.Cells.Item(26, 17).Value = 0

# New trade #26 -> row 27 (momentum)
Set-TextCell $all 27 2 "2026-02-18"
Set-TextCell $all 27 3 "10:29:39"
$all.Cells.Item(27, 1).Value = 26
$all.Cells.Item(27, 4).Value = "momentum"
$all.Cells.Item(27, 5).Value = "UP"
$all.Cells.Item(27, 6).Value = 0.42
$all.Cells.Item(27, 8).Value = "OPEN"
$all.Cells.Item(27, 9).Value = 0
$all.Cells.Item(27, 10).Value = 0
$all.Cells.Item(27, 11).Value = 100
$all.Cells.Item(27, 12).Value = 0
$all.Cells.Item(27, 13).Value = 0
$all.Cells.Item(27, 14).Value = 0.9
$all.Cells.Item(27, 15).Value = "Upward momentum: 1.020% over 5 samples"
$all.Cells.Item(27, 17).Value = 0

# ---------------------------------------------------------------
# "momentum" sheet - only the latest open momentum trade is kept,
# overwritten in place on row 2.
# ---------------------------------------------------------------
$mom = $wb.Worksheets.Item("momentum")
$mom.Cells.Item(2, 1).Value = 26
Set-TextCell $mom 2 3 "10:29:39"
$mom.Cells.Item(2, 5).Value = "UP"
$mom.Cells.Item(2, 6).Value = 0.42
$mom.Cells.Item(2, 15).Value = "Upward momentum: 1.020% over 5 samples"

# ---------------------------------------------------------------
# "MarketMaking" sheet - the 3 most recent open MarketMaking trades,
# overwritten in place on rows 2-4.
# ---------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Cells.Item(2, 1).Value = 23
Set-TextCell $mm 2 3 "10:29:21"
$mm.Cells.Item(2, 6).Value = 0.92
$mm.Cells.Item(2, 15).Value = "Normal spread capture: 408 bps"

$mm.Cells.Item(3, 1).Value = 24
Set-TextCell $mm 3 3 "10:29:26"
$mm.Cells.Item(3, 6).Value = 0.04
$mm.Cells.Item(3, 15).Value = "Normal spread capture: 202 bps"

$mm.Cells.Item(4, 1).Value = 25
Set-TextCell $mm 4 3 "10:29:32"
$mm.Cells.Item(4, 6).Value = 0.388544
$mm.Cells.Item(4, 15).Value = "Normal spread capture: 202 bps"

Write-Output "All Trades / momentum / MarketMaking sheets updated for 2026-02-18 10:30:00"
